$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "81.995.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.90%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "3.186.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'215.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.15%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'624.68"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.47%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "  +19.79%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "  +0.29%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "3.185.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.47%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.591"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.64%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "  +11.62%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "  -3.18%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "3.778.66"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.34%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'31.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.89%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "81.805.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.86%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "3.186.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'3.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +7.41%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'14.05"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.43%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'435.84"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.41%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'8.97"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.76%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'5.14"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.74%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "  +5.73%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'5.34"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +12.09%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "3.349.67"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.49%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'76.75"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.47%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'11.01"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.19%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "  -0.15%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "  +3.77%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'585.03"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +11.22%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'9.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").Value = "'1.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "  +1.40%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "  +8.56%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.27%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.139"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +17.15%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'22.82"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.09%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'6.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.92%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "  +0.44%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("B42").Value = "dogwifhat"
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").Value = "'3.08"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +21.41%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("B43").Value = "Stacks"
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").Value = "'2.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +13.53%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'20.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.85%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'161.03"
$ws.Range("D45").Style = "Normal"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'187.95"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.61%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'44.73"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.34%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.34"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.05%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "  -6.54%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'26.34"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.34%  "
$ws.Range("E51").Style = "Normal"
